# Update the "取得日時" (acquired timestamp) column on the "ランサーズ" sheet
# from the previous run's timestamp to the new append timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$oldValue = "2025-09-17 12:36:24"
$newValue = "2025-09-17 12:45:03"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp = -4162

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Text -eq $oldValue) {
        $cell.Value = $newValue
    }
}
